$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value2 = 9.546140333333334
$ws.Cells.Item(2, 8).Value2 = 28.638421
$ws.Cells.Item(2, 9).Value2 = 0.587227294878132
$ws.Cells.Item(2, 10).Value2 = 0.587227294878132
$ws.Cells.Item(2, 11).Value2 = 2.0
$ws.Cells.Item(2, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(2, 13).Value2 = 0.08324533333333334
$ws.Cells.Item(2, 14).Value2 = 0.249736
$ws.Cells.Item(2, 15).Value2 = 0.05078606388889115
$ws.Cells.Item(2, 16).Value2 = 0.05078606388889115
$ws.Cells.Item(2, 17).Value2 = 0.7946716340951112
$ws.Cells.Item(2, 18).Value2 = 7.152044706856
$ws.Cells.Item(2, 19).Value2 = 0.02982296291498153
$ws.Cells.Item(2, 20).Value2 = 0.02982296291498153

# Row 3
$ws.Cells.Item(3, 7).Value2 = 9.546140333333334
$ws.Cells.Item(3, 8).Value2 = 28.638421
$ws.Cells.Item(3, 9).Value2 = 0.587227294878132
$ws.Cells.Item(3, 10).Value2 = 0.587227294878132
$ws.Cells.Item(3, 15).Value2 = 0.2796082573516313
$ws.Cells.Item(3, 16).Value2 = 0.2796082573516313
$ws.Cells.Item(3, 17).Value2 = 4.375152035058778
$ws.Cells.Item(3, 18).Value2 = 39.376368315529
$ws.Cells.Item(3, 19).Value2 = 0.164193600590187
$ws.Cells.Item(3, 20).Value2 = 0.164193600590187

# Row 4
$ws.Cells.Item(4, 7).Value2 = 9.546140333333334
$ws.Cells.Item(4, 8).Value2 = 28.638421
$ws.Cells.Item(4, 9).Value2 = 0.587227294878132
$ws.Cells.Item(4, 10).Value2 = 0.587227294878132
$ws.Cells.Item(4, 13).Value2 = 1.097575666666667
$ws.Cells.Item(4, 14).Value2 = 3.292727
$ws.Cells.Item(4, 15).Value2 = 0.6696056787594775
$ws.Cells.Item(4, 16).Value2 = 0.6696056787594775
$ws.Cells.Item(4, 17).Value2 = 10.47761134045189
$ws.Cells.Item(4, 18).Value2 = 94.29850206406701
$ws.Cells.Item(4, 19).Value2 = 0.3932107313729634
$ws.Cells.Item(4, 20).Value2 = 0.3932107313729634

# Row 5
$ws.Cells.Item(5, 9).Value2 = 0.2496684258894083
$ws.Cells.Item(5, 10).Value2 = 0.2496684258894083
$ws.Cells.Item(5, 11).Value2 = 2.0
$ws.Cells.Item(5, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(5, 13).Value2 = 0.08324533333333334
$ws.Cells.Item(5, 14).Value2 = 0.249736
$ws.Cells.Item(5, 15).Value2 = 0.05078606388889115
$ws.Cells.Item(5, 16).Value2 = 0.05078606388889115
$ws.Cells.Item(5, 17).Value2 = 0.3378664747262223
$ws.Cells.Item(5, 18).Value2 = 3.040798272536001
$ws.Cells.Item(5, 19).Value2 = 0.01267967662825837
$ws.Cells.Item(5, 20).Value2 = 0.01267967662825837

# Row 6
$ws.Cells.Item(6, 9).Value2 = 0.2496684258894083
$ws.Cells.Item(6, 10).Value2 = 0.2496684258894083
$ws.Cells.Item(6, 15).Value2 = 0.2796082573516313
$ws.Cells.Item(6, 16).Value2 = 0.2796082573516313
$ws.Cells.Item(6, 19).Value2 = 0.06980935347866235
$ws.Cells.Item(6, 20).Value2 = 0.06980935347866234

# Row 7
$ws.Cells.Item(7, 9).Value2 = 0.2496684258894083
$ws.Cells.Item(7, 10).Value2 = 0.2496684258894083
$ws.Cells.Item(7, 13).Value2 = 1.097575666666667
$ws.Cells.Item(7, 14).Value2 = 3.292727
$ws.Cells.Item(7, 15).Value2 = 0.6696056787594775
$ws.Cells.Item(7, 16).Value2 = 0.6696056787594775
$ws.Cells.Item(7, 17).Value2 = 4.454712431230778
$ws.Cells.Item(7, 18).Value2 = 40.09241188107701
$ws.Cells.Item(7, 19).Value2 = 0.1671793957824875
$ws.Cells.Item(7, 20).Value2 = 0.1671793957824875

# Row 8
$ws.Cells.Item(8, 7).Value2 = 2.210442
$ws.Cells.Item(8, 8).Value2 = 6.631326
$ws.Cells.Item(8, 9).Value2 = 0.1359745227725727
$ws.Cells.Item(8, 10).Value2 = 0.1359745227725727
$ws.Cells.Item(8, 11).Value2 = 2.0
$ws.Cells.Item(8, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(8, 13).Value2 = 0.08324533333333334
$ws.Cells.Item(8, 14).Value2 = 0.249736
$ws.Cells.Item(8, 15).Value2 = 0.05078606388889115
$ws.Cells.Item(8, 16).Value2 = 0.05078606388889115
$ws.Cells.Item(8, 17).Value2 = 0.184008981104
$ws.Cells.Item(8, 18).Value2 = 1.656080829936
$ws.Cells.Item(8, 19).Value2 = 0.00690561080078936
$ws.Cells.Item(8, 20).Value2 = 0.00690561080078936

# Row 9
$ws.Cells.Item(9, 7).Value2 = 2.210442
$ws.Cells.Item(9, 8).Value2 = 6.631326
$ws.Cells.Item(9, 9).Value2 = 0.1359745227725727
$ws.Cells.Item(9, 10).Value2 = 0.1359745227725727
$ws.Cells.Item(9, 15).Value2 = 0.2796082573516313
$ws.Cells.Item(9, 16).Value2 = 0.2796082573516313
$ws.Cells.Item(9, 17).Value2 = 1.013081672486
$ws.Cells.Item(9, 18).Value2 = 9.117735052374
$ws.Cells.Item(9, 19).Value2 = 0.03801959935665875
$ws.Cells.Item(9, 20).Value2 = 0.03801959935665874

# Row 10
$ws.Cells.Item(10, 7).Value2 = 2.210442
$ws.Cells.Item(10, 8).Value2 = 6.631326
$ws.Cells.Item(10, 9).Value2 = 0.1359745227725727
$ws.Cells.Item(10, 10).Value2 = 0.1359745227725727
$ws.Cells.Item(10, 13).Value2 = 1.097575666666667
$ws.Cells.Item(10, 14).Value2 = 3.292727
$ws.Cells.Item(10, 15).Value2 = 0.6696056787594775
$ws.Cells.Item(10, 16).Value2 = 0.6696056787594775
$ws.Cells.Item(10, 17).Value2 = 2.426127351778
$ws.Cells.Item(10, 18).Value2 = 21.835146166002
$ws.Cells.Item(10, 19).Value2 = 0.09104931261512454
$ws.Cells.Item(10, 20).Value2 = 0.09104931261512454

# Row 11
$ws.Cells.Item(11, 7).Value2 = 0.4410293333333333
$ws.Cells.Item(11, 8).Value2 = 1.323088
$ws.Cells.Item(11, 9).Value2 = 0.02712975645988715
$ws.Cells.Item(11, 10).Value2 = 0.02712975645988715
$ws.Cells.Item(11, 11).Value2 = 2.0
$ws.Cells.Item(11, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(11, 13).Value2 = 0.08324533333333334
$ws.Cells.Item(11, 14).Value2 = 0.249736
$ws.Cells.Item(11, 15).Value2 = 0.05078606388889115
$ws.Cells.Item(11, 16).Value2 = 0.05078606388889115
$ws.Cells.Item(11, 17).Value2 = 0.03671363386311111
$ws.Cells.Item(11, 18).Value2 = 0.330422704768
$ws.Cells.Item(11, 19).Value2 = 0.001377813544861886
$ws.Cells.Item(11, 20).Value2 = 0.001377813544861886

# Row 12
$ws.Cells.Item(12, 7).Value2 = 0.4410293333333333
$ws.Cells.Item(12, 8).Value2 = 1.323088
$ws.Cells.Item(12, 9).Value2 = 0.02712975645988715
$ws.Cells.Item(12, 10).Value2 = 0.02712975645988715
$ws.Cells.Item(12, 15).Value2 = 0.2796082573516313
$ws.Cells.Item(12, 16).Value2 = 0.2796082573516313
$ws.Cells.Item(12, 17).Value2 = 0.2021309469457778
$ws.Cells.Item(12, 18).Value2 = 1.819178522512
$ws.Cells.Item(12, 19).Value2 = 0.007585703926123208
$ws.Cells.Item(12, 20).Value2 = 0.007585703926123207

# Row 13
$ws.Cells.Item(13, 7).Value2 = 0.4410293333333333
$ws.Cells.Item(13, 8).Value2 = 1.323088
$ws.Cells.Item(13, 9).Value2 = 0.02712975645988715
$ws.Cells.Item(13, 10).Value2 = 0.02712975645988715
$ws.Cells.Item(13, 13).Value2 = 1.097575666666667
$ws.Cells.Item(13, 14).Value2 = 3.292727
$ws.Cells.Item(13, 15).Value2 = 0.6696056787594775
$ws.Cells.Item(13, 16).Value2 = 0.6696056787594775
$ws.Cells.Item(13, 17).Value2 = 0.4840630645528889
$ws.Cells.Item(13, 18).Value2 = 4.356567580976001
$ws.Cells.Item(13, 19).Value2 = 0.01816623898890205
$ws.Cells.Item(13, 20).Value2 = 0.01816623898890206
